# Insert a new weekly price record as row 126 ("Hortaliza, Vega Central
# Mapocho de Santiago - Arveja Verde" sheet), pushing the former rows
# 126-138 down to 127-139 (dimension grows from A1:R138 to A1:R139).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 126 (and everything below it) down by one row.
$ws.Rows(126).Insert()

# Populate the newly inserted row 126 with the new observation.
$ws.Range("A126").Value = 9
$ws.Range("B126").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C126").Value = "Metropolitana"
$ws.Range("D126").Value = 44826
$ws.Range("E126").Value = 13
$ws.Range("F126").Value = 100112022
$ws.Range("G126").Value = "Arveja Verde"
$ws.Range("H126").Value = "Perfection"
$ws.Range("I126").Value = "Primera"
$ws.Range("J126").Value = 45
$ws.Range("K126").Value = 29000
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = 29444
$ws.Range("N126").Value = "`$/malla 25 kilos"
$ws.Range("O126").Value = "Provincia de Limarí"
$ws.Range("P126").Value = 1178
$ws.Range("Q126").Value = 25
$ws.Range("R126").Value = "Hortaliza"
